# "Cleaned up Excel File. V3 BCGW Connection is working"
# Mark every queued AST job in the ast_config table as COMPLETE.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ast_config")

$table = $ws.ListObjects.Item("Table1")
$statusRange = $table.ListColumns.Item("ast_condition").DataBodyRange

foreach ($cell in $statusRange.Cells) {
    $cell.Value = "COMPLETE"
}
